$d = $word.ActiveDocument

function Insert-XmlAfter($range, $bodyFragment) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1) New list-paragraph about panchayat names written in Unicode, inserted
#    right after "Strip new line trails existing in multiple columns."
# ---------------------------------------------------------------------------
$idx = 0
$target = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Strip new line trails existing in multiple columns.*") {
        $target = $idx
    }
}
$stripPara = $d.Paragraphs.Item($target)
$stripPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item($target + 1)

$unicodeFrag = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Certain state</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> other tha</w:t></w:r><w:r><w:t>n Rajasthan had some observations w</w:t></w:r><w:r><w:t>h</w:t></w:r><w:r><w:t>ere panchayat names were written in Unicode. This issue has been addressed</w:t></w:r><w:r><w:t xml:space="preserve"> by the same format in which </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Unicodes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for the panchayat names in Rajasthan was addressed. (In actual process, this issue was identified only during LGD mapping while </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Unicodes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in Rajasthan was identified at an earlier stage and addressed. The same process has been used here.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>)</w:t></w:r></w:p>'
Insert-XmlAfter $newPara.Range $unicodeFrag

# ---------------------------------------------------------------------------
# 2) Drop the stray lastRenderedPageBreak before "Dadra & Nagar Haveli"
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Dadra*Nagar Haveli*") {
        $dadraFrag = '<w:p><w:r><w:t>Dadra &amp; Nagar Haveli</w:t></w:r></w:p>'
        Insert-XmlAfter $p.Range $dadraFrag
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Reflow / de-duplicate the runs in the Madhya Pradesh "BABAI CHICHLI"
#    note and drop the stray _GoBack bookmark inside it.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*BABAI CHICHLI*") {
        $babaiFrag = '<w:p><w:pPr><w:rPr><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve">Missing values in district column for block &#8220;BABAI CHICHLI&#8221;. This belongs to the NARSINGHPUR district which has another set of observations including the above block and I suspect that these are duplicates of the existing ones and also without enough info. </w:t></w:r><w:r><w:rPr><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve">The row with missing observations for the specific </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t>block_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve"> has been filtered out in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t>specific_rows_filter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t>.</w:t></w:r></w:p>'
        Insert-XmlAfter $p.Range $babaiFrag
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Move the lastRenderedPageBreak: it now belongs before "Uttar Pradesh"
#    instead of before "After cleaning".
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Uttar Pradesh*") {
        $upFrag = '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Uttar Pradesh</w:t></w:r></w:p>'
        Insert-XmlAfter $p.Range $upFrag
        break
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "After cleaning*sub directory*") {
        $acFrag = '<w:p><w:r><w:t>After cleaning</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> each state file would be rendered as a new cleaned file bearing its old name and saved in the sub directory &#8220;data/interim/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NREGA_assets</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&#8221;.</w:t></w:r></w:p>'
        Insert-XmlAfter $p.Range $acFrag
        break
    }
}
